$wb = $excel.ActiveWorkbook

# --- "2024" sheet: update Day 1-8 numbers (My 1 / My 2 columns) and fill in Day 9 ---
$ws2024 = $wb.Worksheets.Item("2024")

$ws2024.Range("B2").Value = 202439
$ws2024.Range("C2").Value = 14242

$ws2024.Range("B3").Value = 142709
$ws2024.Range("C3").Value = 34815

$ws2024.Range("B4").Value = 128195
$ws2024.Range("C4").Value = 13479

$ws2024.Range("B5").Value = 101422
$ws2024.Range("C5").Value = 8547

$ws2024.Range("B6").Value = 81853
$ws2024.Range("C6").Value = 10427

$ws2024.Range("B7").Value = 58651
$ws2024.Range("C7").Value = 20749

$ws2024.Range("B8").Value = 57438
$ws2024.Range("C8").Value = 3044

$ws2024.Range("B9").Value = 46509
$ws2024.Range("C9").Value = 2234

# Day 9 (row 10) was previously blank - now filled in. The dependent formulas
# (D10/G10/H10/I10/J10) key off ISBLANK() checks on these cells, so the
# values are written twice to make sure the whole chain recalculates
# consistently once every input cell has a real value.
$ws2024.Range("B10").Value = 11261
$ws2024.Range("C10").Value = 8063
$ws2024.Range("E10").Value = 15866
$ws2024.Range("F10").Value = 11126

$ws2024.Range("B10").Value = 11261
$ws2024.Range("C10").Value = 8063
$ws2024.Range("E10").Value = 15866
$ws2024.Range("F10").Value = 11126

# --- "Overall" sheet: mark 2024 Day 9 (columns AH:AK, row 13) as succeeded ("s") ---
$wsOverall = $wb.Worksheets.Item("Overall")
$wsOverall.Range("AH13").Value = "s"
$wsOverall.Range("AI13").Value = "s"
$wsOverall.Range("AJ13").Value = "s"
$wsOverall.Range("AK13").Value = "s"

# --- Update selections to match the saved view state ---
$wsOverall.Activate()
$wsOverall.Range("AL13").Select()

$ws2024.Activate()
$ws2024.Range("B11").Select()
